# Working on google calendar and gui:
#  - Update the due-date text in D6 (shared string) from
#    "Wed Jun 12 00:00:00 CDT 2019" to "Wed May 08 00:00:00 CDT 2019"
#  - Update the assignment id in B6 from 4523 to 4321
#  - Add a new assignment row (row 7): testing / 6543 / test test test / Sat Jun 08 00:00:00 CDT 2019

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the due date shown in D6
$ws.Range("D6").Value = "Wed May 08 00:00:00 CDT 2019"

# Update the numeric id in B6
$ws.Range("B6").Value = 4321

# Append the new assignment as row 7
$ws.Range("A7").Value = "testing"
$ws.Range("B7").Value = 6543
$ws.Range("C7").Value = "test test test"
$ws.Range("D7").Value = "Sat Jun 08 00:00:00 CDT 2019"
